# Apply league-base update (2024-05-02 20:28) - swap the two
# mismatched match records in each pair of adjacent rows so each
# row's id (col A) lines up with the correct match data (col B:AB).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Serbia Prva Liga")

# Rows 87 and 89: swap match-data columns (B:AB); col A (row id) stays put
$ws.Cells.Item(87, 2).Value = 6989684
$ws.Cells.Item(87, 5).Value = "Radnicki Sremska Mitrovica"
$ws.Cells.Item(87, 6).Value = "FK Tekstilac Odzaci"
$ws.Cells.Item(87, 10).Value = 2.4
$ws.Cells.Item(87, 11).Value = 2.8
$ws.Cells.Item(87, 12).Value = 2.9
$ws.Cells.Item(87, 13).Value = 3.75
$ws.Cells.Item(87, 14).Value = 2.9
$ws.Cells.Item(87, 15).Value = 1.95
$ws.Cells.Item(87, 16).Value = 0.5
$ws.Cells.Item(87, 17).Value = 1.8
$ws.Cells.Item(87, 18).Value = 2
$ws.Cells.Item(87, 19).Value = 2
$ws.Cells.Item(87, 20).Value = 1.85
$ws.Cells.Item(87, 21).Value = 1.95
$ws.Cells.Item(87, 22).Value = 2.75
$ws.Cells.Item(87, 25).Value = 0.8
$ws.Cells.Item(87, 26).Value = -1
$ws.Cells.Item(87, 28).Value = 0.95
$ws.Cells.Item(89, 2).Value = 6989515
$ws.Cells.Item(89, 5).Value = "OFK Vrsac"
$ws.Cells.Item(89, 6).Value = "RFK Novi Sad 1921"
$ws.Cells.Item(89, 10).Value = 1.5
$ws.Cells.Item(89, 11).Value = 3.75
$ws.Cells.Item(89, 12).Value = 5.5
$ws.Cells.Item(89, 13).Value = 1.5
$ws.Cells.Item(89, 14).Value = 3.75
$ws.Cells.Item(89, 15).Value = 6
$ws.Cells.Item(89, 16).Value = -1
$ws.Cells.Item(89, 17).Value = 1.825
$ws.Cells.Item(89, 18).Value = 1.975
$ws.Cells.Item(89, 19).Value = 2.25
$ws.Cells.Item(89, 20).Value = 2
$ws.Cells.Item(89, 21).Value = 1.8
$ws.Cells.Item(89, 22).Value = 0.5
$ws.Cells.Item(89, 25).Value = 0
$ws.Cells.Item(89, 26).Value = 0
$ws.Cells.Item(89, 28).Value = 0.8

# Rows 156 and 157: swap match-data columns (B:AB); col A (row id) stays put
$ws.Cells.Item(156, 2).Value = 6989702
$ws.Cells.Item(156, 5).Value = "FK Mladost Gat Novi Sad"
$ws.Cells.Item(156, 6).Value = "FK Macva Sabac"
$ws.Cells.Item(156, 10).Value = 2.875
$ws.Cells.Item(156, 11).Value = 2.875
$ws.Cells.Item(156, 12).Value = 2.4
$ws.Cells.Item(156, 13).Value = 2.2
$ws.Cells.Item(156, 14).Value = 2.8
$ws.Cells.Item(156, 15).Value = 3.4
$ws.Cells.Item(156, 16).Value = -0.25
$ws.Cells.Item(156, 17).Value = 1.95
$ws.Cells.Item(156, 18).Value = 1.85
$ws.Cells.Item(156, 19).Value = 1.5
$ws.Cells.Item(156, 20).Value = 1.875
$ws.Cells.Item(156, 21).Value = 1.925
$ws.Cells.Item(156, 22).Value = 1.2
$ws.Cells.Item(156, 25).Value = 0.95
$ws.Cells.Item(156, 26).Value = -1
$ws.Cells.Item(156, 28).Value = 0.925
$ws.Cells.Item(157, 2).Value = 6989332
$ws.Cells.Item(157, 5).Value = "OFK Belgrade"
$ws.Cells.Item(157, 6).Value = "FK Dubocica"
$ws.Cells.Item(157, 10).Value = 1.4
$ws.Cells.Item(157, 11).Value = 4
$ws.Cells.Item(157, 12).Value = 7
$ws.Cells.Item(157, 13).Value = 1.285
$ws.Cells.Item(157, 14).Value = 4.333
$ws.Cells.Item(157, 15).Value = 11
$ws.Cells.Item(157, 16).Value = -1.5
$ws.Cells.Item(157, 17).Value = 1.85
$ws.Cells.Item(157, 18).Value = 1.95
$ws.Cells.Item(157, 19).Value = 2.5
$ws.Cells.Item(157, 20).Value = 1.95
$ws.Cells.Item(157, 21).Value = 1.85
$ws.Cells.Item(157, 22).Value = 0.2849999999999999
$ws.Cells.Item(157, 25).Value = -1
$ws.Cells.Item(157, 26).Value = 0.95
$ws.Cells.Item(157, 28).Value = 0.8500000000000001

# Rows 224 and 225: swap match-data columns (B:AB); col A (row id) stays put
$ws.Cells.Item(224, 2).Value = 6989653
$ws.Cells.Item(224, 5).Value = "OFK Belgrade"
$ws.Cells.Item(224, 6).Value = "OFK Vrsac"
$ws.Cells.Item(224, 7).Value = 1
$ws.Cells.Item(224, 8).Value = 0
$ws.Cells.Item(224, 9).Value = "H"
$ws.Cells.Item(224, 10).Value = 1.4
$ws.Cells.Item(224, 11).Value = 4
$ws.Cells.Item(224, 12).Value = 6.5
$ws.Cells.Item(224, 13).Value = 2.6
$ws.Cells.Item(224, 14).Value = 3.2
$ws.Cells.Item(224, 15).Value = 2.375
$ws.Cells.Item(224, 17).Value = 2
$ws.Cells.Item(224, 18).Value = 1.8
$ws.Cells.Item(224, 19).Value = 2.25
$ws.Cells.Item(224, 20).Value = 1.8
$ws.Cells.Item(224, 21).Value = 2
$ws.Cells.Item(224, 22).Value = 1.6
$ws.Cells.Item(224, 24).Value = -1
$ws.Cells.Item(224, 25).Value = 1
$ws.Cells.Item(224, 26).Value = -1
$ws.Cells.Item(224, 27).Value = -1
$ws.Cells.Item(224, 28).Value = 1
$ws.Cells.Item(225, 2).Value = 6989546
$ws.Cells.Item(225, 5).Value = "RFK Novi Sad 1921"
$ws.Cells.Item(225, 6).Value = "Sloboda Uzice"
$ws.Cells.Item(225, 7).Value = 0
$ws.Cells.Item(225, 8).Value = 4
$ws.Cells.Item(225, 9).Value = "A"
$ws.Cells.Item(225, 10).Value = 2.4
$ws.Cells.Item(225, 11).Value = 3
$ws.Cells.Item(225, 12).Value = 2.75
$ws.Cells.Item(225, 13).Value = 2.4
$ws.Cells.Item(225, 14).Value = 3
$ws.Cells.Item(225, 15).Value = 2.75
$ws.Cells.Item(225, 17).Value = 1.775
$ws.Cells.Item(225, 18).Value = 2.025
$ws.Cells.Item(225, 19).Value = 2
$ws.Cells.Item(225, 20).Value = 1.85
$ws.Cells.Item(225, 21).Value = 1.95
$ws.Cells.Item(225, 22).Value = -1
$ws.Cells.Item(225, 24).Value = 1.75
$ws.Cells.Item(225, 25).Value = -1
$ws.Cells.Item(225, 26).Value = 1.025
$ws.Cells.Item(225, 27).Value = 0.8500000000000001
$ws.Cells.Item(225, 28).Value = -1
